$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = "HTML tags should be escaped and displayed as text, not rendered."
$ws.Range("F6").Value = "The text ""SROLLED _TEXT"" will not be displayed and scrolled from right to left"

$ws.Range("G6").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 4
